$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.615.61"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").Value = "1.595.88"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.95"
$ws.Range("E5").Value = "  -0.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.515"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.0618"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.246"
$ws.Range("E9").Value = "  -0.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.47"
$ws.Range("E10").Value = "  -1.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0835"
$ws.Range("E11").Value = "  -0.09%  "
$ws.Range("D12").Value = "1.820.13"
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("D13").Value = "1.596.33"
$ws.Range("E13").Value = "  -0.01%  "
$ws.Range("E14").Value = "  -0.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.522"
$ws.Range("E15").Value = "  -0.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.99"
$ws.Range("E16").Value = "  -0.16%  "
$ws.Range("D17").Value = "26.605.11"
$ws.Range("E17").Value = "  -0.14%  "
$ws.Range("D18").Value = "0.0₃0736"
$ws.Range("E18").Value = "  +0.71%  "
$ws.Range("E19").Value = "  +0.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "208.38"
$ws.Range("E20").Value = "  -0.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.04"
$ws.Range("E21").Value = "  +4.82%  "
$ws.Range("E22").Value = "  +0.23%  "
$ws.Range("E23").Value = "  -0.66%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.88"
$ws.Range("E24").Value = "  -0.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.95"
$ws.Range("E25").Value = "  -1.24%  "
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.15"
$ws.Range("E27").Value = "  -0.43%  "
$ws.Range("E28").Value = "  -0.27%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.25"
$ws.Range("E29").Value = "  -0.59%  "
$ws.Range("E30").Value = "  +0.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.15"
$ws.Range("E31").Value = "  -0.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.22"
$ws.Range("E32").Value = "  -0.81%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.93"
$ws.Range("E33").Value = "  +0.39%  "
$ws.Range("D34").Value = "1.273.28"
$ws.Range("E34").Value = "  -1.54%  "
$ws.Range("E35").Value = "  -8.76%  "
$ws.Range("E36").Value = "  +0.54%  "
$ws.Range("E37").Value = "  -0.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0170"
$ws.Range("E38").Value = "  -0.87%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.838"
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +17.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.47"
$ws.Range("E41").Value = "  +1.26%  "
$ws.Range("E42").Value = "  +0.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.783"
$ws.Range("E43").Value = "  -1.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "64.04"
$ws.Range("E44").Value = "  +0.29%  "
$ws.Range("D45").Value = "1.732.55"
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "90.04"
$ws.Range("E46").Value = "  +0.28%  "
$ws.Range("E47").Value = "  -0.75%  "
$ws.Range("E48").Value = "  +3.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0507"
$ws.Range("E49").Value = "  +0.84%  "
$ws.Range("E50").Value = "  +0.33%  "
$ws.Range("E51").Value = "  -1.48%  "
